# Update the table style applied to the three "Table_0"-styled tables
# (slides 14, 15 and 16) from the default Google-Slides-imported style
# {83FD9BE3-15D6-41B8-AC45-D24EA64D8B20} to the PowerPoint built-in style
# {04999657-C321-4072-890F-FEC90FF4D6A4}.
#
# Table styles cannot be assigned through Table.Style (it is read-only in
# the PowerPoint object model) -- the correct API is Table.ApplyStyle(StyleId).

$p = $ppt.ActivePresentation

$targetStyleId = "{04999657-C321-4072-890F-FEC90FF4D6A4}"
$slideIndexes = 14, 15, 16

foreach ($slideIndex in $slideIndexes) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
